# Horarios Linea 141 - actualizacion de datos (scrap 04:01:01)
$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("LP1912")
$ws2 = $wb.Worksheets.Item("LP1912-215")
$ws3 = $wb.Worksheets.Item("6203-6173")

$updateTime = "04:01:01"

# ---------------------------------------------------------------------------
# Sheet 1: LP1912
# ---------------------------------------------------------------------------
$ws1.Range("A2").Value = "Última actualización: $updateTime"
$ws1.Range("A3").Value = "Total filas: 8"

$ws1Rows = @(
    @($updateTime, "04:01", "81_EL PELIGRO", 0, "LP1912"),
    @($updateTime, "04:46", "215A_EL PATO", 45, "LP1912"),
    @($updateTime, "04:53", "11_ETCHEVERRY", 52, "LP1912"),
    @($updateTime, "05:16", "17_ROMERO", 75, "LP1912"),
    @($updateTime, "05:22", "23_HERNANDEZ", 81, "LP1912"),
    @($updateTime, "05:35", "215B_EL PATO", 94, "LP1912"),
    @($updateTime, "05:41", "14_ABASTO", 100, "LP1912"),
    @($updateTime, "05:46", "15_ABASTO", 105, "LP1912")
)

$r = 6
foreach ($row in $ws1Rows) {
    $ws1.Cells.Item($r, 1).Value = $row[0]
    $ws1.Cells.Item($r, 2).Value = $row[1]
    $ws1.Cells.Item($r, 3).Value = $row[2]
    $ws1.Cells.Item($r, 4).Value = $row[3]
    $ws1.Cells.Item($r, 5).Value = $row[4]
    $r = $r + 1
}

# ---------------------------------------------------------------------------
# Sheet 2: LP1912-215
# ---------------------------------------------------------------------------
$ws2.Range("A2").Value = "Última actualización: $updateTime"
$ws2.Range("A3").Value = "Total filas: 2"

$ws2Rows = @(
    @($updateTime, "04:46", "215A_EL PATO", 45, "LP1912"),
    @($updateTime, "05:35", "215B_EL PATO", 94, "LP1912")
)

$r = 6
foreach ($row in $ws2Rows) {
    $ws2.Cells.Item($r, 1).Value = $row[0]
    $ws2.Cells.Item($r, 2).Value = $row[1]
    $ws2.Cells.Item($r, 3).Value = $row[2]
    $ws2.Cells.Item($r, 4).Value = $row[3]
    $ws2.Cells.Item($r, 5).Value = $row[4]
    $r = $r + 1
}

# ---------------------------------------------------------------------------
# Sheet 3: 6203-6173
# ---------------------------------------------------------------------------
$ws3.Range("A2").Value = "Última actualización: $updateTime"
$ws3.Range("A3").Value = "Total filas: 1"

# Row 5/6 did not exist yet on this sheet - bring over the header formatting
# (bold + border style) and the plain data-row formatting from sheet1 so the
# new rows look like the ones on the other tabs.
$ws1.Range("A5:E5").Copy($ws3.Range("A5:E5"))
$ws1.Range("A6:E6").Copy($ws3.Range("A6:E6"))

$ws3.Cells.Item(6, 1).Value = $updateTime
$ws3.Cells.Item(6, 2).Value = "05:44"
$ws3.Cells.Item(6, 3).Value = "215A_LA PLATA"
$ws3.Cells.Item(6, 4).Value = 103
$ws3.Cells.Item(6, 5).Value = "L6173"
